$d = $word.ActiveDocument
$p1 = $d.Paragraphs(77)
$p3 = $d.Paragraphs(79)
$full = $d.Range($p1.Range.Start, $p3.Range.End)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:tab/><w:t>else:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>$a2 = %string_reg</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:tab/></w:r><w:r><w:t xml:space="preserve">         </w:t></w:r><w:r><w:t>$a3 = null</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$full.InsertXML($xml)

$p2b = $d.Paragraphs(78)
Write-Host "para78 after: [" $p2b.Range.Text "]"
$bstart = $p2b.Range.End - 1 - "$a3 = null".Length
Write-Host "bookmark pos: " $bstart
$bm = $d.Range($bstart, $bstart)
Write-Host "bm text around: [" $bm.Text "]"
$d.Bookmarks.Add("_GoBack", $bm)
